# Scheduled-runner market data refresh: update cached price/profit figures
# across the leve-crafting sheets (ALC, ARM, BSM, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3661.7
$ws.Range("I70").Value = 1194.4117
$ws.Range("J70").Value = 6888.154
$ws.Range("K70").Value = 3583.2351
$ws.Range("L70").Value = 20664.462
$ws.Range("M70").Value = -3313.2351
$ws.Range("N70").Value = -21204.462
$ws.Range("H73").Value = 3661.7
$ws.Range("I73").Value = 1194.4117
$ws.Range("J73").Value = 6888.154
$ws.Range("K73").Value = 3583.2351
$ws.Range("L73").Value = 20664.462
$ws.Range("M73").Value = -2647.2351
$ws.Range("N73").Value = -22536.462
$ws.Range("H74").Value = 84837910
$ws.Range("I74").Value = 84837910
$ws.Range("K74").Value = 84837910
$ws.Range("M74").Value = -84836974
$ws.Range("H76").Value = 3245
$ws.Range("I76").Value = 3245
$ws.Range("K76").Value = 3245
$ws.Range("M76").Value = -2930
$ws.Range("H77").Value = 84837910
$ws.Range("I77").Value = 84837910
$ws.Range("K77").Value = 424189550
$ws.Range("M77").Value = -424184870
$ws.Range("H79").Value = 3245
$ws.Range("I79").Value = 3245
$ws.Range("K79").Value = 3245
$ws.Range("M79").Value = -2153
$ws.Range("H100").Value = 1843.7084
$ws.Range("I100").Value = 1630.5
$ws.Range("J100").Value = 2483.3333
$ws.Range("K100").Value = 1630.5
$ws.Range("L100").Value = 2483.3333
$ws.Range("M100").Value = -1089.5
$ws.Range("N100").Value = -3565.3333
$ws.Range("H106").Value = 2506.7144
$ws.Range("I106").Value = 2591.8462
$ws.Range("K106").Value = 2591.8462
$ws.Range("M106").Value = -1960.8462
$ws.Range("H118").Value = 1743.4286
$ws.Range("I118").Value = 1743.4286
$ws.Range("K118").Value = 5230.2858
$ws.Range("M118").Value = -3573.2858
$ws.Range("H127").Value = 1123
$ws.Range("I127").Value = 1123
$ws.Range("K127").Value = 3369
$ws.Range("M127").Value = 1591
$ws.Range("H137").Value = 1925.7241
$ws.Range("I137").Value = 1437.3478
$ws.Range("J137").Value = 3797.8333
$ws.Range("K137").Value = 4312.0434
$ws.Range("L137").Value = 11393.4999
$ws.Range("M137").Value = -1762.0434
$ws.Range("N137").Value = -16493.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 885.6539
$ws.Range("I97").Value = 1010.75
$ws.Range("K97").Value = 1010.75
$ws.Range("M97").Value = -514.75
$ws.Range("H132").Value = 5885478
$ws.Range("I132").Value = 7144937.5
$ws.Range("K132").Value = 21434812.5
$ws.Range("M132").Value = -21432282.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1975.5625
$ws.Range("I99").Value = 1649.875
$ws.Range("J99").Value = 2301.25
$ws.Range("K99").Value = 1649.875
$ws.Range("L99").Value = 2301.25
$ws.Range("M99").Value = -151.875
$ws.Range("N99").Value = -5297.25
$ws.Range("H134").Value = 29417162
$ws.Range("I134").Value = 31255484
$ws.Range("K134").Value = 93766452
$ws.Range("M134").Value = -93763917

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 708147.5
$ws.Range("I4").Value = 1001209
$ws.Range("K4").Value = 3003627
$ws.Range("M4").Value = -3003515
$ws.Range("H23").Value = 43.25
$ws.Range("J23").Value = 47.5
$ws.Range("L23").Value = 142.5
$ws.Range("N23").Value = -612.5
$ws.Range("H34").Value = 451.5
$ws.Range("I34").Value = 458
$ws.Range("J34").Value = 440.66666
$ws.Range("K34").Value = 1374
$ws.Range("L34").Value = 1321.99998
$ws.Range("M34").Value = -1290
$ws.Range("N34").Value = -1489.99998
$ws.Range("H39").Value = 1133.1666
$ws.Range("I39").Value = 700.25
$ws.Range("J39").Value = 1999
$ws.Range("K39").Value = 2100.75
$ws.Range("L39").Value = 5997
$ws.Range("M39").Value = -1806.75
$ws.Range("N39").Value = -6585
$ws.Range("H55").Value = 1190
$ws.Range("I55").Value = 1475
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 4425
$ws.Range("L55").Value = 3000
$ws.Range("M55").Value = -4248
$ws.Range("N55").Value = -3354
$ws.Range("H68").Value = 1613.5714
$ws.Range("J68").Value = 2500
$ws.Range("L68").Value = 7500
$ws.Range("N68").Value = -9122
$ws.Range("H70").Value = 10128.866
$ws.Range("I70").Value = 5093.4
$ws.Range("K70").Value = 15280.2
$ws.Range("M70").Value = -14965.2
$ws.Range("H71").Value = 1613.5714
$ws.Range("J71").Value = 2500
$ws.Range("L71").Value = 22500
$ws.Range("N71").Value = -30612
$ws.Range("H73").Value = 10128.866
$ws.Range("I73").Value = 5093.4
$ws.Range("K73").Value = 15280.2
$ws.Range("M73").Value = -14188.2
$ws.Range("H80").Value = 6493.5
$ws.Range("I80").Value = 7241.75
$ws.Range("J80").Value = 4997
$ws.Range("K80").Value = 21725.25
$ws.Range("L80").Value = 14991
$ws.Range("M80").Value = -20789.25
$ws.Range("N80").Value = -16863
$ws.Range("H83").Value = 6493.5
$ws.Range("I83").Value = 7241.75
$ws.Range("J83").Value = 4997
$ws.Range("K83").Value = 65175.75
$ws.Range("L83").Value = 44973
$ws.Range("M83").Value = -60495.75
$ws.Range("N83").Value = -54333
$ws.Range("H129").Value = 3806.647
$ws.Range("J129").Value = 4966.5
$ws.Range("L129").Value = 14899.5
$ws.Range("N129").Value = -24899.5
$ws.Range("H140").Value = 2887.0908
$ws.Range("I140").Value = 2887.0908
$ws.Range("K140").Value = 8661.2724
$ws.Range("M140").Value = -3481.2724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 19999
$ws.Range("J6").Value = 19999
$ws.Range("L6").Value = 19999
$ws.Range("N6").Value = -20225
$ws.Range("H11").Value = 27250
$ws.Range("I11").Value = 40000
$ws.Range("K11").Value = 40000
$ws.Range("M11").Value = -39861
$ws.Range("H16").Value = 19999
$ws.Range("J16").Value = 19999
$ws.Range("L16").Value = 19999
$ws.Range("N16").Value = -20499
$ws.Range("H126").Value = 3550
$ws.Range("H132").Value = 7355628.5
$ws.Range("I132").Value = 8931214
$ws.Range("J132").Value = 2899
$ws.Range("K132").Value = 26793642
$ws.Range("L132").Value = 8697
$ws.Range("M132").Value = -26791112
$ws.Range("N132").Value = -13757

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1865.6957
$ws.Range("I82").Value = 1928.1333
$ws.Range("J82").Value = 1748.625
$ws.Range("K82").Value = 1928.1333
$ws.Range("L82").Value = 1748.625
$ws.Range("M82").Value = -1567.1333
$ws.Range("N82").Value = -2470.625
$ws.Range("H85").Value = 1865.6957
$ws.Range("I85").Value = 1928.1333
$ws.Range("J85").Value = 1748.625
$ws.Range("K85").Value = 1928.1333
$ws.Range("L85").Value = 1748.625
$ws.Range("M85").Value = -680.1333
$ws.Range("N85").Value = -4244.625
$ws.Range("H122").Value = 14715.429
$ws.Range("I122").Value = 14329.333
$ws.Range("K122").Value = 42987.999
$ws.Range("M122").Value = -40537.999
$ws.Range("H132").Value = 24003606
$ws.Range("I132").Value = 28238826
$ws.Range("J132").Value = 4031.6667
$ws.Range("K132").Value = 84716478
$ws.Range("L132").Value = 12095.0001
$ws.Range("M132").Value = -84713948
$ws.Range("N132").Value = -17155.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3414.037
$ws.Range("I96").Value = 1456.1428
$ws.Range("K96").Value = 1456.1428
$ws.Range("M96").Value = -83.14280000000008
$ws.Range("H100").Value = 2158.6
$ws.Range("I100").Value = 2065.111
$ws.Range("K100").Value = 4130.222
$ws.Range("M100").Value = -3589.222
$ws.Range("H107").Value = 443.6842
$ws.Range("I107").Value = 370.625
$ws.Range("K107").Value = 1111.875
$ws.Range("M107").Value = 808.125
$ws.Range("H122").Value = 1744.909
$ws.Range("I122").Value = 1719.4
$ws.Range("K122").Value = 5158.200000000001
$ws.Range("M122").Value = -2708.200000000001
$ws.Range("H132").Value = 17872596
$ws.Range("J132").Value = 41200
$ws.Range("L132").Value = 123600
$ws.Range("N132").Value = -128660
